$d = $word.ActiveDocument

# Locate the "WORD:" list paragraph (the PDF entry must be inserted right before it,
# mirroring the existing "HTML:" list item).
$wordParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^WORD:") {
        $wordParaIndex = $i
        break
    }
}

$wordPara = $d.Paragraphs.Item($wordParaIndex)

# Insert a new (empty) paragraph before the "WORD:" item; Word copies the
# paragraph formatting (Compact style + the same numbered-list membership)
# from the following paragraph, exactly like the other list entries.
$wordPara.Range.InsertParagraphBefore()

$pdfPara = $d.Paragraphs.Item($wordParaIndex)
$pdfRange = $pdfPara.Range
$pdfRange.Collapse(1)

# Type the label and the URL as plain text first ...
$pdfRange.InsertAfter("PDF: http://sbennett1990.github.io/Resume/resume.pdf")

# ... then turn the URL portion into a real hyperlink, the same way the
# "HTML:" and "WORD:" entries do it.
$paraRange = $pdfPara.Range
$urlStart = $paraRange.Start + 5
$urlEnd = $paraRange.End - 1
$urlRange = $d.Range($urlStart, $urlEnd)
$d.Hyperlinks.Add($urlRange, "http://sbennett1990.github.io/Resume/resume.pdf", [Type]::Missing, [Type]::Missing, [Type]::Missing) | Out-Null

# Add the trailing line break that separates this list item from the next,
# matching the "HTML:" entry above it.
$endRange = $pdfPara.Range
$endRange.Collapse(0)
$endRange.InsertBreak(6)
